# Fix issue mapping Medication ab77bff992840669d45583ae812eee5175aff7fe
#
# 1. Update the "Date" metadata value on the "Metadata" sheet.
# 2. On the "Mapping Table 0" sheet: clear the Source text ("Messages/...Forme")
#    from A4 and A5 (duplicate rows that should no longer carry that source),
#    and append a new mapping row (row 6) with Relationship "related-to" and
#    Target "Medication.form.coding.code", matching the formatting of row 5.

$wb = $excel.ActiveWorkbook

# --- Update metadata Date value --------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-11-17T10:38:58+00:00"

# --- Update the mapping table -----------------------------------------------
$wsMap = $wb.Worksheets.Item("Mapping Table 0")

# Clear the Source column text on rows 4 and 5 (keep formatting/style)
$wsMap.Range("A4").Value = ""
$wsMap.Range("A5").Value = ""

# Create row 6 by copying the formatting of row 5, then set its values
$wsMap.Range("A5:E5").Copy()
$wsMap.Range("A6:E6").PasteSpecial(-4122)

$wsMap.Range("C6").Value = "related-to"
$wsMap.Range("D6").Value = "Medication.form.coding.code"
